$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-03-11 03:04:27"
$wsZhCn.Range("G2").Value = "2016-03-11 03:04:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-03-11 03:04:35"
$wsDeDe.Range("G2").Value = "2016-03-11 03:05:18"
